$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 14:46"

# Update country rows: label (col A) + 7 numeric stats (cols B:H)
# (row, country, totalCases, newCases, activeCases, recovered, criticalCases, deathsToday, deaths)
$rowsData = @(
    (4, "Estados Unidos", 4568484, 447, 2245521, 2169115, 0, 8, 153848),
    (6, "India", 1590297, 5913, 1023811, 531438, 0, 45, 35048),
    (40, "Kuwait", 66529, 626, 57330, 8754, 0, 1, 445),
    (41, "Republica Dominicana", 66182, 0, 33947, 31112, 0, 0, 1123),
    (44, "Paises Bajos", 53963, 342, 0, 0, 0, 0, 6147),
    (46, "Portugal", 50868, 255, 36140, 13001, 0, 2, 1727),
    (64, "Uzbekistan", 23078, 493, 13432, 9512, 0, 3, 134),
    (67, "Nepal", 19547, 274, 14248, 5247, 0, 3, 52),
    (72, "Australia", 16303, 721, 10619, 5495, 0, 13, 189),
    (78, "Dinamarca", 13725, 91, 12526, 584, 0, 1, 615),
    (79, "Estado de Palestina", 11548, 264, 5016, 6451, 0, 1, 81),
    (80, "Sudan", 11496, 0, 6001, 4770, 0, 0, 725),
    (81, "Bosnia y Herzegovina", 11444, 317, 5586, 5530, 0, 12, 328),
    (82, "Bulgaria", 11155, 0, 5971, 4816, 0, 0, 368),
    (84, "Republica de Macedonia", 10617, 114, 6020, 4117, 0, 4, 480),
    (85, "Senegal", 10106, 145, 6725, 3177, 0, 4, 204),
    (97, "Zambia", 5555, 306, 3289, 2117, 0, 3, 149),
    (100, "Croacia", 5071, 78, 4178, 749, 0, 3, 144),
    (131, "Islandia", 1872, 11, 1823, 39, 0, 0, 10),
    (144, "Uganda", 1147, 7, 1028, 117, 0, 0, 2),
    (146, "Burkina Faso", 1106, 1, 935, 118, 0, 0, 53),
    (165, "Gambia", 403, 77, 66, 329, 0, 0, 8),
    (166, "Guyana", 398, 0, 185, 193, 0, 0, 20),
    (167, "Burundi", 387, 0, 304, 82, 0, 0, 1),
    (168, "Comoras", 378, 0, 330, 41, 0, 0, 7),
    (169, "Birmania", 353, 2, 294, 53, 0, 0, 6),
    (170, "Mauricio", 344, 0, 332, 2, 0, 0, 10),
    (171, "Isla de Man", 336, 0, 312, 0, 0, 0, 24),
    (179, "Gibraltar", 187, 1, 180, 7, 0, 0, 0)
)

foreach ($r in $rowsData) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
}